$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1316.4445
$ws.Range("J19").Value = 906.8570999999999
$ws.Range("L19").Value = 906.8570999999999
$ws.Range("N19").Value = -1256.8571
$ws.Range("H33").Value = 614.2105
$ws.Range("I33").Value = 569.3570999999999
$ws.Range("K33").Value = 569.3570999999999
$ws.Range("M33").Value = -340.3570999999999
$ws.Range("H58").Value = 1292.8235
$ws.Range("I58").Value = 99
$ws.Range("J58").Value = 1548.6428
$ws.Range("K58").Value = 297
$ws.Range("L58").Value = 4645.928400000001
$ws.Range("M58").Value = -147
$ws.Range("N58").Value = -4945.928400000001
$ws.Range("H76").Value = 6286.4
$ws.Range("I76").Value = 6033.3335
$ws.Range("K76").Value = 6033.3335
$ws.Range("M76").Value = -5718.3335
$ws.Range("H79").Value = 6286.4
$ws.Range("I79").Value = 6033.3335
$ws.Range("K79").Value = 6033.3335
$ws.Range("M79").Value = -4941.3335
$ws.Range("H98").Value = 992.5
$ws.Range("I98").Value = 992.5
$ws.Range("K98").Value = 992.5
$ws.Range("M98").Value = 505.5
$ws.Range("H99").Value = 254.42857
$ws.Range("J99").Value = 181
$ws.Range("L99").Value = 543
$ws.Range("N99").Value = -3539
$ws.Range("H101").Value = 875
$ws.Range("I101").Value = 750
$ws.Range("J101").Value = 1000
$ws.Range("K101").Value = 2250
$ws.Range("L101").Value = 3000
$ws.Range("M101").Value = -628
$ws.Range("N101").Value = -6244
$ws.Range("H106").Value = 6543
$ws.Range("I106").Value = 6406.2856
$ws.Range("K106").Value = 6406.2856
$ws.Range("M106").Value = -5775.2856
$ws.Range("H118").Value = 1484.75
$ws.Range("J118").Value = 209
$ws.Range("L118").Value = 627
$ws.Range("N118").Value = -3941
$ws.Range("H122").Value = 992.5
$ws.Range("I122").Value = 992.5
$ws.Range("K122").Value = 2977.5
$ws.Range("M122").Value = -527.5
$ws.Range("H127").Value = 4774.7334
$ws.Range("J127").Value = 8832.833000000001
$ws.Range("L127").Value = 26498.499
$ws.Range("N127").Value = -36418.499
$ws.Range("H135").Value = 1522.1052
$ws.Range("I135").Value = 1420.0625
$ws.Range("J135").Value = 2066.3333
$ws.Range("K135").Value = 12780.5625
$ws.Range("L135").Value = 18596.9997
$ws.Range("M135").Value = -10245.5625
$ws.Range("N135").Value = -23666.9997
$ws.Range("H137").Value = 5035.5405
$ws.Range("I137").Value = 5908.524
$ws.Range("J137").Value = 3889.75
$ws.Range("K137").Value = 17725.572
$ws.Range("L137").Value = 11669.25
$ws.Range("M137").Value = -15175.572
$ws.Range("N137").Value = -16769.25
$ws.Range("H138").Value = 3141.5454
$ws.Range("J138").Value = 4189.7
$ws.Range("L138").Value = 12569.1
$ws.Range("N138").Value = -22849.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5775
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 5775
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 5775
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -6199
$ws.Range("H74").Value = 1381.2142
$ws.Range("I74").Value = 1426.091
$ws.Range("J74").Value = 1216.6666
$ws.Range("K74").Value = 1426.091
$ws.Range("L74").Value = 1216.6666
$ws.Range("M74").Value = -552.0909999999999
$ws.Range("N74").Value = -2964.6666
$ws.Range("H77").Value = 1381.2142
$ws.Range("I77").Value = 1426.091
$ws.Range("J77").Value = 1216.6666
$ws.Range("K77").Value = 7130.455
$ws.Range("L77").Value = 6083.333000000001
$ws.Range("M77").Value = -2762.455
$ws.Range("N77").Value = -14819.333
$ws.Range("H132").Value = 25643350
$ws.Range("I132").Value = 35716556
$ws.Range("J132").Value = 2461.818
$ws.Range("K132").Value = 107149668
$ws.Range("L132").Value = 7385.454000000001
$ws.Range("M132").Value = -107147138
$ws.Range("N132").Value = -12445.454
$ws.Range("H136").Value = 5775
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5775
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 17325
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -22425

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1932.8077
$ws.Range("I31").Value = 1655.6316
$ws.Range("K31").Value = 1655.6316
$ws.Range("M31").Value = -1360.6316
$ws.Range("H34").Value = 1932.8077
$ws.Range("I34").Value = 1655.6316
$ws.Range("K34").Value = 1655.6316
$ws.Range("M34").Value = -1453.6316
$ws.Range("H58").Value = 1892.3334
$ws.Range("I58").Value = 1427.7142
$ws.Range("K58").Value = 1427.7142
$ws.Range("M58").Value = -1224.7142
$ws.Range("H132").Value = 2564.7742
$ws.Range("I132").Value = 1960.6
$ws.Range("J132").Value = 5082.1665
$ws.Range("K132").Value = 5881.799999999999
$ws.Range("L132").Value = 15246.4995
$ws.Range("M132").Value = -3351.799999999999
$ws.Range("N132").Value = -20306.4995
$ws.Range("H134").Value = 2054.7568
$ws.Range("I134").Value = 1641.8387
$ws.Range("K134").Value = 4925.5161
$ws.Range("M134").Value = -2390.5161
$ws.Range("H136").Value = 1892.3334
$ws.Range("I136").Value = 1427.7142
$ws.Range("K136").Value = 4283.142599999999
$ws.Range("M136").Value = -1733.142599999999
$ws.Range("H138").Value = 40943.05
$ws.Range("J138").Value = 39366.668
$ws.Range("L138").Value = 39366.668
$ws.Range("N138").Value = -49646.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9973.615
$ws.Range("I3").Value = 7165.7
$ws.Range("J3").Value = 19333.334
$ws.Range("K3").Value = 21497.1
$ws.Range("L3").Value = 58000.00199999999
$ws.Range("M3").Value = -21385.1
$ws.Range("N3").Value = -58224.00199999999
$ws.Range("H5").Value = 981.8889
$ws.Range("I5").Value = 878.3333
$ws.Range("J5").Value = 1499.6666
$ws.Range("K5").Value = 2634.9999
$ws.Range("L5").Value = 4498.9998
$ws.Range("M5").Value = -2522.9999
$ws.Range("N5").Value = -4722.9998
$ws.Range("H14").Value = 503.33334
$ws.Range("I14").Value = 503.33334
$ws.Range("K14").Value = 1510.00002
$ws.Range("M14").Value = -1337.00002
$ws.Range("H68").Value = 15153806
$ws.Range("I68").Value = 600
$ws.Range("J68").Value = 16669127
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 50007381
$ws.Range("M68").Value = -989
$ws.Range("N68").Value = -50009003
$ws.Range("H71").Value = 15153806
$ws.Range("I71").Value = 600
$ws.Range("J71").Value = 16669127
$ws.Range("K71").Value = 5400
$ws.Range("L71").Value = 150022143
$ws.Range("M71").Value = -1344
$ws.Range("N71").Value = -150030255
$ws.Range("H107").Value = 12823006
$ws.Range("I107").Value = 2658.8333
$ws.Range("K107").Value = 7976.499899999999
$ws.Range("M107").Value = -6056.499899999999
$ws.Range("H134").Value = 3513.4167
$ws.Range("I134").Value = 2105.5454
$ws.Range("K134").Value = 6316.6362
$ws.Range("M134").Value = -1246.6362
$ws.Range("H135").Value = 981.8889
$ws.Range("I135").Value = 878.3333
$ws.Range("J135").Value = 1499.6666
$ws.Range("K135").Value = 7904.9997
$ws.Range("L135").Value = 13496.9994
$ws.Range("M135").Value = -5369.9997
$ws.Range("N135").Value = -18566.9994
$ws.Range("H139").Value = 3121.2856
$ws.Range("I139").Value = 3121.2856
$ws.Range("K139").Value = 9363.856800000001
$ws.Range("M139").Value = -4223.856800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8040.0835
$ws.Range("I113").Value = 5022.625
$ws.Range("K113").Value = 5022.625
$ws.Range("M113").Value = -2852.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 330.84375
$ws.Range("I55").Value = 407.5
$ws.Range("K55").Value = 407.5
$ws.Range("M55").Value = -234.5
$ws.Range("H61").Value = 1023.2381
$ws.Range("I61").Value = 1044.4
$ws.Range("J61").Value = 600
$ws.Range("K61").Value = 1044.4
$ws.Range("L61").Value = 600
$ws.Range("M61").Value = -842.4000000000001
$ws.Range("N61").Value = -1004
$ws.Range("H68").Value = 2675
$ws.Range("J68").Value = 2675
$ws.Range("L68").Value = 2675
$ws.Range("N68").Value = -4173
$ws.Range("H71").Value = 2675
$ws.Range("J71").Value = 2675
$ws.Range("L71").Value = 13375
$ws.Range("N71").Value = -20863
$ws.Range("H113").Value = 1023.2381
$ws.Range("I113").Value = 1044.4
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 1044.4
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 1125.6
$ws.Range("N113").Value = -4940
$ws.Range("H132").Value = 3080.9048
$ws.Range("I132").Value = 2480.8823
$ws.Range("K132").Value = 7442.646900000001
$ws.Range("M132").Value = -4912.646900000001
$ws.Range("H136").Value = 3759.625
$ws.Range("I136").Value = 3087.375
$ws.Range("K136").Value = 9262.125
$ws.Range("M136").Value = -6712.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 428.7143
$ws.Range("I113").Value = 280.9
$ws.Range("K113").Value = 842.6999999999999
$ws.Range("M113").Value = 1327.3
$ws.Range("H132").Value = 3079.1936
$ws.Range("I132").Value = 2683.9092
$ws.Range("J132").Value = 4045.4443
$ws.Range("K132").Value = 8051.7276
$ws.Range("L132").Value = 12136.3329
$ws.Range("M132").Value = -5521.7276
$ws.Range("N132").Value = -17196.3329
$ws.Range("H136").Value = 2103.5386
$ws.Range("I136").Value = 2103.5386
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6310.6158
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3760.6158
$ws.Range("N136").ClearContents()
